$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph (the Jekyll/Jupiter footer block
# that should be removed, together with the blank paragraph that precedes
# it and the "(c) 2020 ..." paragraph that follows it).
$rng = $d.Content
$find = $rng.Find
$find.Text = "Ver no Jupiter Salvar em pdf Salvar em docx"
$find.Execute() | Out-Null

if ($find.Found -eq $true) {
    $jupiterIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            $jupiterIndex = $i
            break
        }
    }

    if ($jupiterIndex -gt 0) {
        $startPara = $d.Paragraphs.Item($jupiterIndex - 1)
        $endPara = $d.Paragraphs.Item($jupiterIndex + 1)

        $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $deleteRange.Delete()
    }
}
